# Refresh the crypto price/volume snapshot table on Sheet1 (coinranking.com export).
# Columns: A=index, B=Coin, C=Link, D=Price, E=Volume(1h), F=Data, G=Hora.
# D/E hold numeric-looking text ("314.79", "2.18%") that must stay literal text
# (matches the source file's inlineStr cells), so those assignments use a leading
# apostrophe -- the same trick Excel's UI uses to force text entry instead of
# auto-converting to a Number/Percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'314.79"
$ws.Range("E2").Value = "'2.18%"

$ws.Range("D3").Value = "'39.24"
$ws.Range("E3").Value = "'-1.84%"

$ws.Range("D4").Value = "'5.144"
$ws.Range("E4").Value = "'-0.04%"

$ws.Range("D5").Value = "'0.08168"
$ws.Range("E5").Value = "'0.31%"

$ws.Range("D6").Value = "'1.989"
$ws.Range("E6").Value = "'2.36%"

$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.378"
$ws.Range("E7").Value = "'3.23%"

$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'8.334"
$ws.Range("E8").Value = "'2.20%"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9360"
$ws.Range("E9").Value = "'0.68%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1311"
$ws.Range("E10").Value = "'-8.82%"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1969"
$ws.Range("E11").Value = "'2.47%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08969"
$ws.Range("E12").Value = "'-1.82%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03525"
$ws.Range("E13").Value = "'0.19%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09735"
$ws.Range("E14").Value = "'-0.48%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001407"
$ws.Range("E15").Value = "'0.37%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006513"
$ws.Range("E16").Value = "'10.81%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.632"
$ws.Range("E17").Value = "'-7.31%"

$ws.Range("D18").Value = "'3.125"
$ws.Range("E18").Value = "'-7.50%"

$ws.Range("D19").Value = "'0.3471"
$ws.Range("E19").Value = "'1.22%"

$ws.Range("D20").Value = "'0.1317"
$ws.Range("E20").Value = "'0.43%"

$ws.Range("D21").Value = "'4.983"
$ws.Range("E21").Value = "'7.45%"

$ws.Range("D22").Value = "'0.2489"
$ws.Range("E22").Value = "'2.67%"

$ws.Range("D23").Value = "'0.04375"
$ws.Range("E23").Value = "'0.03%"

$ws.Range("D24").Value = "'0.001243"
$ws.Range("E24").Value = "'1.08%"

$ws.Range("D25").Value = "'0.004763"
$ws.Range("E25").Value = "'8.88%"

$ws.Range("D26").Value = "'0.0003891"
$ws.Range("E26").Value = "'198.98%"

$ws.Range("E27").Value = "'-7.64%"

$ws.Range("D39").Value = "'0.02240"
$ws.Range("E39").Value = "'9.44%"

$ws.Range("D40").Value = "'0.05198"
$ws.Range("E40").Value = "'2.59%"

$ws.Range("E41").Value = "'4.00%"

$ws.Range("D42").Value = "'0.01029"
$ws.Range("E42").Value = "'4.33%"

$ws.Range("D43").Value = "'0.1397"
$ws.Range("E43").Value = "'2.28%"

$ws.Range("D44").Value = "'0.002101"
$ws.Range("E44").Value = "'-1.50%"

$ws.Range("D45").Value = "'0.008857"
$ws.Range("E45").Value = "'-5.52%"

$ws.Range("D46").Value = "'0.00006820"
$ws.Range("E46").Value = "'6.93%"

$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.06%"

$ws.Range("D48").Value = "'0.003007"
$ws.Range("E48").Value = "'10.77%"

$ws.Range("E49").Value = "'30.00%"

$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.06%"

$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.06%"
